# Insert a new data row at row 8 (shifting existing rows 8-80 down to 9-81)
# and populate it with the new price-report record for this date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("8:8").Insert()

$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 44537
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 100112032
$ws.Range("G8").Value = "Zapallo italiano"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 170
$ws.Range("K8").Value = 5500
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5765
$ws.Range("N8").Value = "$/caja 60 unidades"
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 96
$ws.Range("Q8").Value = 60
$ws.Range("R8").Value = "Hortaliza"
